$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.226.20"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "2.999.70"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.49%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "3.002.34"
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.35%  "
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("E13").Value = "  -3.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.63%  "
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "3.496.15"
$ws.Range("E16").Value = "  -2.08%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "62.172.81"
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("D19").Value = "3.001.34"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.682"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.55%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("E32").Value = "  -1.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  -4.44%  "
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("D36").Value = "0.0₃0782"
$ws.Range("E36").Value = "  -4.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.73%  "
$ws.Range("E38").Value = "  -4.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "415.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.86%  "
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.274"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.95%  "
$ws.Range("D45").Value = "2.761.65"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0349"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.108"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.60%  "
